$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Vtn","Plaur","ECs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,40.92730033333333,122.781901,0.3921621516522625,0.3921621516522625,270.0027617753436,2430.024855978093,0.04743372348041777,0.04743372348041776),
  @("ECs","Vtn","Plaur","FAPs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,7.598228333333334,22.794685,0.07280562235174674,0.07280562235174673,50.12650768291167,451.138569146205,0.008806157717929673,0.00880615771792967),
  @("ECs","Vtn","Plaur","M2",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,49.40125466666666,148.203764,0.4733589112063359,0.4733589112063358,325.9065486003613,2933.158937403251,0.05725482585852041,0.0572548258585204),
  @("ECs","Vtn","Plaur","sCs",2,0.6666666666666666,6.597131,19.791393,0.1209543635982448,0.1209543635982448,3,1,6.436425,19.309275,0.06167331478965488,0.06167331478965488,42.461938896675,382.157450070075,0.007459656541376924,0.007459656541376923),
  @("FAPs","Vtn","Plaur","ECs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,40.92730033333333,122.781901,0.3921621516522625,0.3921621516522625,861.6495479830577,7754.84593184752,0.1513734382837992,0.1513734382837992),
  @("FAPs","Vtn","Plaur","FAPs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,7.598228333333334,22.794685,0.07280562235174674,0.07280562235174673,159.9668181279111,1439.7013631512,0.02810275630971167,0.02810275630971167),
  @("FAPs","Vtn","Plaur","M2",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,49.40125466666666,148.203764,0.4733589112063359,0.4733589112063358,1040.053177381476,9360.478596433279,0.1827151488987024,0.1827151488987024),
  @("FAPs","Vtn","Plaur","sCs",3,1,21.05317333333333,63.15952,0.3859970617919927,0.3859970617919927,3,1,6.436425,19.309275,0.06167331478965488,0.06167331478965488,135.507171172,1219.564540548,0.02380571829977943,0.02380571829977943),
  @("sCs","Vtn","Plaur","ECs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,40.92730033333333,122.781901,0.3921621516522625,0.3921621516522625,1100.61739712187,9905.556574096832,0.1933549898880456,0.1933549898880455),
  @("sCs","Vtn","Plaur","FAPs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,7.598228333333334,22.794685,0.07280562235174674,0.07280562235174673,204.3316373877689,1838.98473648992,0.03589670832410539,0.03589670832410539),
  @("sCs","Vtn","Plaur","M2",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,49.40125466666666,148.203764,0.4733589112063359,0.4733589112063358,1328.499067442716,11956.49160698445,0.233388936449113,0.233388936449113),
  @("sCs","Vtn","Plaur","sCs",3,1,26.89201066666667,80.676032,0.4930485746097625,0.4930485746097625,3,1,6.436425,19.309275,0.06167331478965488,0.06167331478965488,173.0884097552,1557.7956877968,0.03040793994849852,0.03040793994849852)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($c = 1; $c -le $row.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}

